# Apply odds updates to Sheet1 of the FlashScore workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 changes
$ws.Range("Q3").Value = 2.25
$ws.Range("R3").Value = 1.62

# Row 4 changes
$ws.Range("G4").Value = 2.8
$ws.Range("H4").Value = 3.1
$ws.Range("I4").Value = 2.45
$ws.Range("J4").Value = 3.25
$ws.Range("L4").Value = 3
$ws.Range("O4").Value = 1.3
$ws.Range("X4").Value = 15.5
$ws.Range("Y4").Value = 10
$ws.Range("AA4").Value = 23
$ws.Range("AH4").Value = 7.9
$ws.Range("AI4").Value = 12
$ws.Range("AK4").Value = 27
$ws.Range("AO4").Value = 14.5
$ws.Range("AT4").Value = 2.62
$ws.Range("AU4").Value = 6.4
$ws.Range("AW4").Value = 4.45

$wb.Save()
